# "Add Init; change cable-15"
# Updates the cable-15 worksheet: the "Quantity" (Qty for 1 cable) value in
# C23 goes from 8 to 20. The dependent formulas in D27:D30 (=$C$23*Cxx)
# recalculate automatically. Also updates the sheet's active selection to
# C24 (from G8), matching the last cell the author left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantity: 8 -> 20 (drives D27:D30 via =$C$23*Cxx formulas)
$ws.Range("C23").Value = 20

# Move the active selection to C24 (was G8), and scroll so row 5 is at the
# top of the view.
$ws.Range("C24").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
